# Update cryptocurrency price/volume data per latest refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.844.48"
$ws.Range("E2").Value = "  -2.81%  "

$ws.Range("D3").Value = "'1.792.52"
$ws.Range("E3").Value = "  -0.61%  "

$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").Value = "'316.43"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.32%  "

$ws.Range("D7").Value = "'0.5306"
$ws.Range("E7").Value = "  -0.31%  "

$ws.Range("D8").Value = "'0.3855"
$ws.Range("E8").Value = "  +2.16%  "

$ws.Range("D9").Value = "'0.07442"
$ws.Range("E9").Value = "  -0.69%  "

$ws.Range("D10").Value = "'41.44"
$ws.Range("E10").Value = "  -2.55%  "

$ws.Range("D11").Value = "'1.086"
$ws.Range("E11").Value = "  -2.17%  "

$ws.Range("D12").Value = "'1.0000"
$ws.Range("E12").Value = "  -0.35%  "

$ws.Range("D13").Value = "'6.191"
$ws.Range("E13").Value = "  +0.71%  "

$ws.Range("D14").Value = "'7.462"
$ws.Range("E14").Value = "  +1.67%  "

$ws.Range("D15").Value = "'20.33"
$ws.Range("E15").Value = "  -1.59%  "

$ws.Range("D16").Value = "'1.791.03"
$ws.Range("E16").Value = "  -0.74%  "

$ws.Range("D17").Value = "'88.23"
$ws.Range("E17").Value = "  -2.25%  "

$ws.Range("D18").Value = "'0.00001059"
$ws.Range("E18").Value = "  -0.47%  "

$ws.Range("D19").Value = "'0.06518"
$ws.Range("E19").Value = "  +1.00%  "

$ws.Range("D20").Value = "'0.9997"
$ws.Range("E20").Value = "  -0.28%  "

$ws.Range("D21").Value = "'17.25"
$ws.Range("E21").Value = "  +0.27%  "

$ws.Range("D22").Value = "'5.956"
$ws.Range("E22").Value = "  +0.85%  "

$ws.Range("D23").Value = "'27.891.57"
$ws.Range("E23").Value = "  -2.70%  "

$ws.Range("E24").Value = "  +0.89%  "

$ws.Range("D25").Value = "'2.092"
$ws.Range("E25").Value = "  -0.23%  "

$ws.Range("D26").Value = "'157.17"
$ws.Range("E26").Value = "  -1.60%  "

$ws.Range("D27").Value = "'20.16"
$ws.Range("E27").Value = "  -1.45%  "

$ws.Range("D28").Value = "'1.996.84"
$ws.Range("E28").Value = "  -0.74%  "

$ws.Range("D29").Value = "'2.326"
$ws.Range("E29").Value = "  -1.46%  "

$ws.Range("D30").Value = "'121.45"
$ws.Range("E30").Value = "  -0.96%  "

$ws.Range("D31").Value = "'0.1093"
$ws.Range("E31").Value = "  +4.40%  "

$ws.Range("D32").Value = "'1.103"
$ws.Range("E32").Value = "  +0.64%  "

$ws.Range("D33").Value = "'3.649"
$ws.Range("E33").Value = "  -1.01%  "

$ws.Range("D34").Value = "'5.506"
$ws.Range("E34").Value = "  -2.28%  "

$ws.Range("D35").Value = "'0.06925"
$ws.Range("E35").Value = "  +8.01%  "

$ws.Range("D36").Value = "'0.2203"
$ws.Range("E36").Value = "  -2.04%  "

$ws.Range("D37").Value = "'0.02269"
$ws.Range("E37").Value = "  -1.34%  "

$ws.Range("D38").Value = "'5.043"
$ws.Range("E38").Value = "  +0.40%  "

$ws.Range("D39").Value = "'11.41"
$ws.Range("E39").Value = "  +1.28%  "

$ws.Range("D40").Value = "'8.380"
$ws.Range("E40").Value = "  -4.90%  "

$ws.Range("D41").Value = "'0.6104"
$ws.Range("E41").Value = "  -1.77%  "

$ws.Range("D42").Value = "'1.173"
$ws.Range("E42").Value = "  -4.85%  "

$ws.Range("E43").Value = "  +0.33%  "

$ws.Range("D44").Value = "'13.32"
$ws.Range("E44").Value = "  +0.64%  "

$ws.Range("D45").Value = "'3.678"
$ws.Range("E45").Value = "  -0.34%  "

$ws.Range("D46").Value = "'0.5696"
$ws.Range("E46").Value = "  -2.56%  "

$ws.Range("D47").Value = "'124.87"
$ws.Range("E47").Value = "  -0.69%  "

$ws.Range("D48").Value = "'1.910"
$ws.Range("E48").Value = "  -1.56%  "

$ws.Range("D49").Value = "'1.170"
$ws.Range("E49").Value = "  +2.10%  "

$ws.Range("D50").Value = "'0.06793"
$ws.Range("E50").Value = "  -1.38%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'71.32"
$ws.Range("E51").Value = "  -1.17%  "
